$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.906.71'
$ws.Range("E2").Value = '  +0.15%  '
$ws.Range("D3").Value = '1.889.03'
$ws.Range("E3").Value = '  -0.12%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.7669'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.07%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '242.81'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.60%  '
$ws.Range("E7").Value = '  -0.08%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3135'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.10%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '25.71'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.58%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07168'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.78%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08520'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +4.70%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.7637'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.23%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.375'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.74%  '
$ws.Range("D14").Value = '1.844.51'
$ws.Range("E14").Value = '  -0.78%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '93.92'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.47%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.156'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.74%  '
$ws.Range("D17").Value = '29.732.36'
$ws.Range("E17").Value = '  -0.41%  '
$ws.Range("E18").Value = '  -1.06%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '244.34'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.21%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000007798'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.87%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.9991'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.04%  '
$ws.Range("D22").Value = '2.122.96'
$ws.Range("E22").Value = '  +0.09%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.032'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.48%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.9999'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.04%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1622'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +3.66%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.406'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.16%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '162.32'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.52%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.76'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.01%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.034'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.12%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.499'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.87%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.536'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.64%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.494'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.14%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.105'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.32%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05428'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.80%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.241'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.65%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7429'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.06%  '
$ws.Range("E37").Value = '  +0.38%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.695'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.79%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01947'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.50%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.780'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.30%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.4469'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.41%  '
$ws.Range("D42").Value = '1.099.33'
$ws.Range("E42").Value = '  -3.93%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.080'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.99%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '72.92'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.14%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.8527'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.04%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.9998'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.04%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '103.04'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.06%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.870'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.55%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.676'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.36%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '3.008'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -3.57%  '
$ws.Range("D51").Value = '2.016.04'
$ws.Range("E51").Value = '  -0.65%  '
